$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function ReplaceParaXML($index, [string[]]$innerXmlList) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $xml = ($innerXmlList -join "")
    $r.InsertXML($xml) | Out-Null
}

# NOTE: process paragraphs from last to first so inserting brand-new
# paragraphs doesn't shift the Paragraphs.Item(N) index of paragraphs
# that haven't been processed yet.

# ---- Old paragraph 9 -> New paragraph 10 + New paragraph 11 (new) ----
$p10 = "<w:p $wns><w:r><w:t xml:space=`"preserve`">10. ¿Qué información debe incluirse en el email que se le envía al transportista? ¿y en la notificación </w:t></w:r><w:proofErr $wns w:type=`"spellStart`"/><w:r><w:t>push</w:t></w:r><w:proofErr $wns w:type=`"spellEnd`"/><w:r><w:t>?</w:t></w:r></w:p>"
$p11 = "<w:p $wns><w:r><w:t>1</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:t xml:space=`"preserve`">. </w:t></w:r><w:r><w:t>¿</w:t></w:r><w:r><w:t>Prefiere visualizar los pedidos y cotizaciones en pantallas diferentes?</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r></w:p>"
ReplaceParaXML 9 @($p10, $p11)

# ---- Old paragraph 8 -> New paragraph 9 ----
$p9 = "<w:p $wns><w:r><w:t>9. ¿Qué pasarelas de pago estarán disponibles? las mismas se deben mostrar en un menú desplegable, una lista con casilleros de selección, ¿o de qué manera?</w:t></w:r></w:p>"
ReplaceParaXML 8 @($p9)

# ---- Old paragraph 7 -> New paragraph 8 ----
$p8 = "<w:p $wns><w:r><w:t>8. Para seleccionar la forma de pago, las mismas se deben mostrar en un menú desplegable, una lista con casilleros de selección, ¿o de qué manera?</w:t></w:r></w:p>"
ReplaceParaXML 7 @($p8)

# ---- Old paragraph 6 -> New paragraph 7 ----
$p7 = "<w:p $wns><w:r><w:t>7. Una vez aceptada una cotización, ¿se deben mostrar los datos de contacto de los usuarios?</w:t></w:r></w:p>"
ReplaceParaXML 6 @($p7)

# ---- Old paragraph 5 -> New paragraph 6 ----
$p6 = "<w:p $wns><w:r><w:t>6. Una vez aceptada una cotización, ¿se puede cancelar?</w:t></w:r></w:p>"
ReplaceParaXML 5 @($p6)

# ---- Old paragraph 4 -> New paragraph 5 ----
$p5 = "<w:p $wns><w:r><w:t>5. ¿El alcance de esta primera entrega comienza visualizando todas las cotizaciones existentes?</w:t></w:r></w:p>"
ReplaceParaXML 4 @($p5)

# ---- Old paragraph 3 -> New paragraph 4 ----
$p4 = "<w:p $wns><w:r><w:t>4. ¿Se debe informar que el estado del pedido de envío se cambió a `"Confirmado`"?</w:t></w:r></w:p>"
ReplaceParaXML 3 @($p4)

# ---- Old paragraph 2 -> New paragraph 2 + New paragraph 3 (new) ----
$p2 = "<w:p $wns><w:r><w:t xml:space=`"preserve`">2. ¿Cuál es el formato que prefiere utilizar para visualizar la fecha? Por </w:t></w:r><w:r><w:t>ejemplo,</w:t></w:r><w:r><w:t xml:space=`"preserve`"> el formato </w:t></w:r><w:proofErr $wns w:type=`"spellStart`"/><w:r><w:t>dd</w:t></w:r><w:proofErr $wns w:type=`"spellEnd`"/><w:r><w:t>/mm/</w:t></w:r><w:proofErr $wns w:type=`"spellStart`"/><w:r><w:t>aaaa</w:t></w:r><w:proofErr $wns w:type=`"spellEnd`"/><w:r><w:t>.</w:t></w:r></w:p>"
$p3 = "<w:p $wns><w:r><w:t xml:space=`"preserve`">3. Para ingresar los datos tipo fecha, ¿prefiere utilizar un calendario en el que pueda seleccionar la fecha u otra manera? </w:t></w:r></w:p>"
ReplaceParaXML 2 @($p2, $p3)

# ---- Old paragraph 1 -> New paragraph 1 ----
$p1 = "<w:p $wns><w:r><w:t>1. ¿Qué paleta de colores le gustaría utilizar?</w:t></w:r></w:p>"
ReplaceParaXML 1 @($p1)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
